# Update the product catalog on Hoja1 (sheet1):
#  - fix the existing "Piñata de Bluei" row's description/image text
#  - append three new product rows (Bingo, Plim Plim, Stich)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2: correct existing "Piñata de Bluei" entry ---
$ws.Cells.Item(2, 3).Value = "piñata color azul, 80 de alto por 40 de ancho"
$ws.Cells.Item(2, 5).Value = "Piñatas/bluei.jpg"

# --- Row 3: new "Piñata de Bingo" entry ---
$ws.Cells.Item(3, 1).Value = "Piñata de Bingo"
$ws.Cells.Item(3, 2).Value = 100
$ws.Cells.Item(3, 3).Value = "piñata color naranja, 80 de alto por 40 de ancho"
$ws.Cells.Item(3, 4).Value = "Piñatas"
$ws.Cells.Item(3, 5).Value = "Piñatas/bingo.jpg"

# --- Row 4: new "Piñata de Plim Plim" entry (no category set) ---
$ws.Cells.Item(4, 1).Value = "Piñata de Plim Plim"
$ws.Cells.Item(4, 2).Value = 100
$ws.Cells.Item(4, 3).Value = "piñata color roja con blanco, 80 cm de altura por 30 de ancho."
$ws.Cells.Item(4, 5).Value = "Piñatas/plimplim.jpg"

# --- Row 5: new "Piñata de Stich" entry (no category set) ---
$ws.Cells.Item(5, 1).Value = "Piñata de Stich"
$ws.Cells.Item(5, 2).Value = 100
$ws.Cells.Item(5, 3).Value = "piñata color azul, 80 cm de altura por 30 de ancho."
$ws.Cells.Item(5, 5).Value = "Piñatas/stich.jpg"
